# The deck's live theme (ppt/theme/theme1.xml, "Integral") is switched back
# to the default Office color palette ("Office Theme") by re-writing every
# entry of the master's 12-slot theme color scheme (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink). Font scheme and format scheme already matched
# the Office defaults, so only the colors need to change.
#
# RGB() below packs standard 0xRRGGBB hex into PowerPoint's BGR-ordered
# OLE_COLOR integer (0x00BBGGRR), matching ColorFormat.RGB semantics.
function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$colors.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$colors.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
